$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 950
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 450
$ws.Range("B5").Value = 152
$ws.Range("B6").Value = 85
$ws.Range("B8").Value = 65
